# Adds the "Python Code" test-case rows to Sheet2, and creates a new
# Sheet4 holding the same rows (mirrors the upstream commit that added
# these tree-module test cases to the workbook).

$wb = $excel.ActiveWorkbook

# ---- 1. Extend Sheet2 with 3 new rows (5-7) ----
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("A5").Value = "Pyhton Code:"

$ws2.Range("A6").Value = "Valid Python Code"
$ws2.Range("B6").Value = 'print("This is a valid python code")'

$ws2.Range("A7").Value = "Incali Python Code"
$ws2.Range("B7").Value = 'p("Invalid Code")'

# ---- 2. Add a new worksheet "Sheet4" after the last sheet ----
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "Sheet4"

$ws4.Range("A1").Value = "Pyhton Code:"

$ws4.Range("A2").Value = "Valid Python Code"
$ws4.Range("B2").Value = 'print("This is a valid python code")'

$ws4.Range("A3").Value = "Incali Python Code"
$ws4.Range("B3").Value = 'p("Invalid Code")'
